$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Gurkeerat Singh Mann" followed by a trailing NBSP (U+00A0), matching the
# existing F2:F6 values exactly.
$batsman = "Gurkeerat Singh Mann" + [char]0x00A0

# Rows 7-11 duplicate rows 2-6 (venue, date, result, ownTeam, oppTeam,
# batsman, totalRuns, totalBalls, total4s, total6s, sr).
$rows = @(
    @(" Dubai (DSC)", " October 25 2020", "Super Kings won by 8 wickets (with 8 balls remaining)", "Royal Challengers Bangalore", "Chennai Super Kings", $batsman, "2", "2", "0", "0", "100.00"),
    @(" Abu Dhabi", " October 28 2020", "Mumbai won by 5 wickets (with 5 balls remaining)", "Royal Challengers Bangalore", "Mumbai Indians", $batsman, "14", "11", "2", "0", "127.27"),
    @(" Sharjah", " October 31 2020", "Sunrisers won by 5 wickets (with 35 balls remaining)", "Royal Challengers Bangalore", "Sunrisers Hyderabad", $batsman, "15", "24", "1", "0", "62.50"),
    @(" Dubai (DSC)", " October 17 2020", "RCB won by 7 wickets (with 2 balls remaining)", "Royal Challengers Bangalore", "Rajasthan Royals", $batsman, "19", "17", "1", "0", "111.76"),
    @(" Abu Dhabi", " October 21 2020", "RCB won by 8 wickets (with 39 balls remaining)", "Royal Challengers Bangalore", "Kolkata Knight Riders", $batsman, "21", "26", "4", "0", "80.76")
)

$startRow = 7
for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $startRow + $i
    $rowData = $rows[$i]
    for ($c = 1; $c -le $rowData.Length; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        # Force text storage (matching the source sheet, where even numeric
        # looking values like "2" or "100.00" are stored as text), then
        # restore the default "Normal" style so no stray formatting is left
        # behind on the cell.
        $cell.NumberFormat = "@"
        $cell.Value = $rowData[$c - 1]
        $cell.Style = "Normal"
    }
}
